$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells to remain text (matches original inline string
# cells which are NOT real numbers, e.g. "27.493.47" / "1.842.17" use dots as
# thousands separators and would otherwise be auto-converted by Excel).
$priceCells = @("D2","D3","D5","D7","D8","D9","D10","D11","D13","D14","D15","D16","D18","D22","D23","D24","D26","D27","D28","D29","D30","D31","D32","D34","D36","D38","D39","D40","D41","D43","D44","D45","D46","D47","D48","D49","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.493.47"
$ws.Range("E2").Value = "  -2.08%  "
$ws.Range("D3").Value = "1.842.17"
$ws.Range("E3").Value = "  -2.81%  "
$ws.Range("E4").Value = "  -1.13%  "
$ws.Range("D5").Value = "332.92"
$ws.Range("E5").Value = "  -1.43%  "
$ws.Range("E6").Value = "  -1.05%  "
$ws.Range("D7").Value = "0.4618"
$ws.Range("E7").Value = "  -2.93%  "
$ws.Range("D8").Value = "0.3831"
$ws.Range("E8").Value = "  -3.62%  "
$ws.Range("D9").Value = "46.26"
$ws.Range("E9").Value = "  -2.83%  "
$ws.Range("D10").Value = "0.07879"
$ws.Range("E10").Value = "  -2.34%  "
$ws.Range("D11").Value = "0.9784"
$ws.Range("E11").Value = "  -4.59%  "
$ws.Range("E12").Value = "  -4.12%  "
$ws.Range("D13").Value = "1.844.53"
$ws.Range("E13").Value = "  -2.69%  "
$ws.Range("D14").Value = "5.896"
$ws.Range("E14").Value = "  -2.50%  "
$ws.Range("D15").Value = "7.029"
$ws.Range("E15").Value = "  -3.17%  "
$ws.Range("D16").Value = "1.005"
$ws.Range("E16").Value = "  -1.08%  "
$ws.Range("E17").Value = "  -1.08%  "
$ws.Range("D18").Value = "0.06625"
$ws.Range("E18").Value = "  -2.17%  "
$ws.Range("E19").Value = "  -2.49%  "
$ws.Range("E20").Value = "  -1.37%  "
$ws.Range("E21").Value = "  -0.92%  "
$ws.Range("D22").Value = "27.497.99"
$ws.Range("E22").Value = "  -1.99%  "
$ws.Range("D23").Value = "5.344"
$ws.Range("E23").Value = "  -3.88%  "
$ws.Range("D24").Value = "10.87"
$ws.Range("E24").Value = "  -1.83%  "
$ws.Range("E25").Value = "  -2.52%  "
$ws.Range("D26").Value = "157.11"
$ws.Range("E26").Value = "  -2.49%  "
$ws.Range("D27").Value = "19.38"
$ws.Range("E27").Value = "  -3.80%  "
$ws.Range("D28").Value = "2.073"
$ws.Range("E28").Value = "  -2.45%  "
$ws.Range("D29").Value = "5.353"
$ws.Range("E29").Value = "  -3.96%  "
$ws.Range("D30").Value = "118.93"
$ws.Range("E30").Value = "  -2.68%  "
$ws.Range("D31").Value = "0.9570"
$ws.Range("E31").Value = "  -2.99%  "
$ws.Range("D32").Value = "0.09316"
$ws.Range("E33").Value = "  -2.08%  "
$ws.Range("D34").Value = "5.242"
$ws.Range("E34").Value = "  -2.63%  "
$ws.Range("E35").Value = "  -4.10%  "
$ws.Range("D36").Value = "0.05942"
$ws.Range("E36").Value = "  -2.91%  "
$ws.Range("E37").Value = "  -3.07%  "
$ws.Range("D38").Value = "8.102"
$ws.Range("E38").Value = "  -2.01%  "
$ws.Range("D39").Value = "1.159"
$ws.Range("E39").Value = "  -4.16%  "
$ws.Range("D40").Value = "0.5834"
$ws.Range("E40").Value = "  -3.04%  "
$ws.Range("D41").Value = "0.1843"
$ws.Range("E41").Value = "  -3.24%  "
$ws.Range("E42").Value = "  -3.06%  "
$ws.Range("D43").Value = "1.248"
$ws.Range("E43").Value = "  -2.46%  "
$ws.Range("D44").Value = "0.5505"
$ws.Range("E44").Value = "  -3.42%  "
$ws.Range("D45").Value = "11.98"
$ws.Range("E45").Value = "  -2.35%  "
$ws.Range("D46").Value = "1.870"
$ws.Range("E46").Value = "  -3.83%  "
$ws.Range("D47").Value = "0.06667"
$ws.Range("E47").Value = "  -2.55%  "
$ws.Range("D48").Value = "110.26"
$ws.Range("E48").Value = "  -2.23%  "
$ws.Range("D49").Value = "1.043"
$ws.Range("E49").Value = "  -3.34%  "
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("D51").Value = "69.72"
$ws.Range("E51").Value = "  -2.64%  "
